$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (row 9) to the daily routine log.
$ws.Range("A9").Value = "12/12/2019 Data flow,front screens for admin  ,faculty"
$ws.Range("C9").Value = "dataflow diagram for the project was drawn"
$ws.Range("B9").Value = "dataflow,front screens for modules were drwan."

# Move the active selection to B9, matching the saved selection state.
$ws.Range("B9").Select()
